# Generate Report for Handoff
#
# A new source file (ff4e45d6-f564-42f7-b61f-da98018a26ad.md) has shown up
# "Ready for handoff" and needs a row inserted right before the trailing
# ".localization-config" bookkeeping row on every sheet (Overview, zh-cn,
# de-de). Inserting the row pushes ".localization-config" down by one.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # OLE BGR encoding of RGB FF6495ED (the workbook's HyperLink font color)

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name) / B (zh-cn) / C (de-de)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Push the ".localization-config" row from row 8 down to row 9, values +
# styles included.
$wsOverview.Rows.Item(8).Insert()

# Row 9 now holds the old ".localization-config" row's values/style; its
# hyperlink target needs to be (re)anchored at A9.
$wsOverview.Hyperlinks.Add($wsOverview.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/6167ab42a53096a0714091546d107070aef7e174/.localization-config", "", "", ".localization-config") | Out-Null

# Row 8: the new file, "Ready for handoff" in both locale columns.
$wsOverview.Range("B8").Value = "Ready for handoff"
$wsOverview.Range("C8").Value = "Ready for handoff"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/6167ab42a53096a0714091546d107070aef7e174/e2e/ff4e45d6-f564-42f7-b61f-da98018a26ad.md", "", "", "ff4e45d6-f564-42f7-b61f-da98018a26ad.md") | Out-Null
Style-AsHyperlink $wsOverview.Range("A8")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(8).Insert()

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/6167ab42a53096a0714091546d107070aef7e174/.localization-config", "", "", ".localization-config") | Out-Null

$wsZh.Range("B8").Value = "Ready for handoff"
$wsZh.Range("D8").Value = "2016-02-25 05:35:21"
$wsZh.Range("G8").Value = "0001-01-01 00:00:00"
$wsZh.Range("H8").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/6167ab42a53096a0714091546d107070aef7e174/e2e/ff4e45d6-f564-42f7-b61f-da98018a26ad.md", "", "", "ff4e45d6-f564-42f7-b61f-da98018a26ad.md") | Out-Null
Style-AsHyperlink $wsZh.Range("A8")

$wsZh.Hyperlinks.Add($wsZh.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0bc30c1088d5bdb51b63dc922db27a8f5f8b34e7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ff4e45d6-f564-42f7-b61f-da98018a26ad.0bc30c1088d5bdb51b63dc922db27a8f5f8b34e7.zh-cn.xlf", "", "", "ff4e45d6-f564-42f7-b61f-da98018a26ad.0bc30c1088d5bdb51b63dc922db27a8f5f8b34e7.zh-cn.xlf") | Out-Null
Style-AsHyperlink $wsZh.Range("C8")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(8).Insert()

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/6167ab42a53096a0714091546d107070aef7e174/.localization-config", "", "", ".localization-config") | Out-Null

$wsDe.Range("B8").Value = "Ready for handoff"
$wsDe.Range("D8").Value = "2016-02-25 05:35:33"
$wsDe.Range("G8").Value = "0001-01-01 00:00:00"
$wsDe.Range("H8").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/6167ab42a53096a0714091546d107070aef7e174/e2e/ff4e45d6-f564-42f7-b61f-da98018a26ad.md", "", "", "ff4e45d6-f564-42f7-b61f-da98018a26ad.md") | Out-Null
Style-AsHyperlink $wsDe.Range("A8")

$wsDe.Hyperlinks.Add($wsDe.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0bc30c1088d5bdb51b63dc922db27a8f5f8b34e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ff4e45d6-f564-42f7-b61f-da98018a26ad.0bc30c1088d5bdb51b63dc922db27a8f5f8b34e7.de-de.xlf", "", "", "ff4e45d6-f564-42f7-b61f-da98018a26ad.0bc30c1088d5bdb51b63dc922db27a8f5f8b34e7.de-de.xlf") | Out-Null
Style-AsHyperlink $wsDe.Range("C8")

Write-Host "Report generated for handoff."
